$wb = $excel.ActiveWorkbook

# --- Update the LeaveType sheet: rename headers, restyle ---
$leaveType = $wb.Worksheets.Item("LeaveType")
$leaveType.Range("A1").Value = "LeaveName"
$leaveType.Range("B1").Value = "LeaveAbbrevation"
$leaveType.Range("A1:B1").Interior.Color = 65535

# Reset the selection on LeaveType to B1 (it will no longer be the active tab)
$leaveType.Range("B1").Select() | Out-Null

# --- Add the new EmployeeType sheet right after LeaveType ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $leaveType)
$newSheet.Name = "EmployeeType"

$newSheet.Columns.Item(1).ColumnWidth = 32.71

$newSheet.Range("A1").Value = "EmpCategoryName "
$newSheet.Range("A2").Value = "Permanent"
$newSheet.Range("A3").Value = "Temporary"
$newSheet.Range("A1").Interior.Color = 65535

# EmployeeType becomes the active tab/selected sheet, with C6 selected
$newSheet.Range("C6").Select() | Out-Null
